$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 420
$range = $ws.Range("C2:C" + $lastRow)
$range.Value = 45190
